$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 385.53845
$ws.Range("J41").Value = 603.75
$ws.Range("L41").Value = 603.75
$ws.Range("N41").Value = -1483.75
$ws.Range("H101").Value = 2560
$ws.Range("I101").Value = 2700
$ws.Range("K101").Value = 8100
$ws.Range("M101").Value = -6478
$ws.Range("H103").Value = 3400
$ws.Range("H137").Value = 1795.1818
$ws.Range("I137").Value = 1392.4286
$ws.Range("K137").Value = 4177.2858
$ws.Range("M137").Value = -1627.2858
$ws.Range("H138").Value = 5253.609
$ws.Range("J138").Value = 5418.6924
$ws.Range("L138").Value = 16256.0772
$ws.Range("N138").Value = -26536.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388
$ws.Range("H14").Value = 365
$ws.Range("I14").Value = 365
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 365
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -190
$ws.Range("H32").Value = 3225.0286
$ws.Range("I32").Value = 2073.7273
$ws.Range("K32").Value = 2073.7273
$ws.Range("M32").Value = -1786.7273
$ws.Range("H61").Value = 15664.333
$ws.Range("I61").Value = 19996.5
$ws.Range("K61").Value = 19996.5
$ws.Range("M61").Value = -19784.5
$ws.Range("H74").Value = 1470.2858
$ws.Range("I74").Value = 1414.4546
$ws.Range("K74").Value = 1414.4546
$ws.Range("M74").Value = -540.4546
$ws.Range("H77").Value = 1470.2858
$ws.Range("I77").Value = 1414.4546
$ws.Range("K77").Value = 7072.273
$ws.Range("M77").Value = -2704.273
$ws.Range("H122").Value = 2673
$ws.Range("I122").Value = 2673
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8019
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5569
$ws.Range("H136").Value = 15664.333
$ws.Range("I136").Value = 19996.5
$ws.Range("K136").Value = 59989.5
$ws.Range("M136").Value = -57439.5
$ws.Range("N5").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("H22").Value = 658.4
$ws.Range("I22").Value = 730.8333
$ws.Range("J22").Value = 549.75
$ws.Range("K22").Value = 730.8333
$ws.Range("L22").Value = 549.75
$ws.Range("M22").Value = -557.8333
$ws.Range("N22").Value = -895.75
$ws.Range("H134").Value = 3368
$ws.Range("I134").Value = 1739
$ws.Range("K134").Value = 5217
$ws.Range("M134").Value = -2682
$ws.Range("N4").ClearContents()
$ws.Range("M15").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1193.6666
$ws.Range("I7").Value = 1173.8
$ws.Range("J7").Value = 1211.7273
$ws.Range("K7").Value = 1173.8
$ws.Range("L7").Value = 1211.7273
$ws.Range("M7").Value = -1060.8
$ws.Range("N7").Value = -1437.7273
$ws.Range("H15").Value = 9649.5
$ws.Range("I15").Value = 300
$ws.Range("J15").Value = 18999
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 18999
$ws.Range("M15").Value = -130
$ws.Range("N15").Value = -19339
$ws.Range("H19").Value = 2068
$ws.Range("I19").Value = 75.666664
$ws.Range("J19").Value = 19999
$ws.Range("K19").Value = 75.666664
$ws.Range("L19").Value = 19999
$ws.Range("M19").Value = 94.333336
$ws.Range("N19").Value = -20339
$ws.Range("H22").Value = 1800
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1800
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1800
$ws.Range("N22").Value = -2500
$ws.Range("H24").Value = 2068
$ws.Range("I24").Value = 75.666664
$ws.Range("J24").Value = 19999
$ws.Range("K24").Value = 75.666664
$ws.Range("L24").Value = 19999
$ws.Range("M24").Value = 94.333336
$ws.Range("N24").Value = -20339
$ws.Range("H25").Value = 472.5
$ws.Range("I25").Value = 363.33334
$ws.Range("K25").Value = 363.33334
$ws.Range("M25").Value = -189.33334
$ws.Range("H43").Value = 29749.834
$ws.Range("J43").Value = 29749.834
$ws.Range("L43").Value = 29749.834
$ws.Range("N43").Value = -30117.834
$ws.Range("H101").Value = 29749.834
$ws.Range("J101").Value = 29749.834
$ws.Range("L101").Value = 29749.834
$ws.Range("N101").Value = -36239.834
$ws.Range("H105").Value = 2049.3
$ws.Range("I105").Value = 1498.5
$ws.Range("K105").Value = 1498.5
$ws.Range("M105").Value = 248.5
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 715
$ws.Range("I44").Value = 930
$ws.Range("K44").Value = 2790
$ws.Range("M44").Value = -2392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 346.69232
$ws.Range("I2").Value = 26
$ws.Range("K2").Value = 26
$ws.Range("M2").Value = 87
$ws.Range("H3").Value = 835789
$ws.Range("I3").Value = 2500235
$ws.Range("J3").Value = 3566
$ws.Range("K3").Value = 2500235
$ws.Range("L3").Value = 3566
$ws.Range("M3").Value = -2500119
$ws.Range("N3").Value = -3798
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H43").Value = 6666.5
$ws.Range("I43").Value = 6666.5
$ws.Range("K43").Value = 6666.5
$ws.Range("M43").Value = -6515.5
$ws.Range("H102").Value = 998.25
$ws.Range("I102").Value = 998.25
$ws.Range("K102").Value = 998.25
$ws.Range("M102").Value = 623.75
$ws.Range("H122").Value = 334666
$ws.Range("I122").Value = 500999.5
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 1502998.5
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -1500548.5
$ws.Range("N122").Value = -10897
$ws.Range("M21").ClearContents()
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4144.6665
$ws.Range("I40").Value = 3973.6
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3973.6
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3837.6
$ws.Range("N40").Value = -5272
$ws.Range("H122").Value = 7409.091
$ws.Range("I122").Value = 5680
$ws.Range("J122").Value = 8850
$ws.Range("K122").Value = 17040
$ws.Range("L122").Value = 26550
$ws.Range("M122").Value = -14590
$ws.Range("N122").Value = -31450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1500
$ws.Range("I26").Value = 1500
$ws.Range("K26").Value = 1500
$ws.Range("M26").Value = -1207
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H100").Value = 11616695
$ws.Range("I100").Value = 17424542
$ws.Range("K100").Value = 34849084
$ws.Range("M100").Value = -34848543
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("M122").ClearContents()
